$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 16 (pushing the old rows 16-32 down to 18-34),
# mirroring the two new price records added to the weekly Fruta/hortaliza sheet.
$ws.Range("A16:A17").EntireRow.Insert()

# New row 16: Maracuya "Primera" record dated 2021-08-27 (serial 44435)
$ws.Range("A16").Value = 9
$ws.Range("B16").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C16").Value = "Metropolitana"
$ws.Range("D16").Value = 44435
$ws.Range("E16").Value = 13
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100108
$ws.Range("H16").Value = "Tropicales y subtropicales"
$ws.Range("I16").Value = 100108003
$ws.Range("J16").Value = "Maracuyá"
$ws.Range("K16").Value = "Sin especificar"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 30
$ws.Range("N16").Value = 32000
$ws.Range("O16").Value = 32000
$ws.Range("P16").Value = 32000
$ws.Range("Q16").Value = "$/caja 18 kilos"
$ws.Range("R16").Value = "Perú"
$ws.Range("S16").Value = 1778
$ws.Range("T16").Value = 18

# New row 17: Maracuya "Primera" record dated 2021-08-23 (serial 44431)
$ws.Range("A17").Value = 9
$ws.Range("B17").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C17").Value = "Metropolitana"
$ws.Range("D17").Value = 44431
$ws.Range("E17").Value = 13
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100108
$ws.Range("H17").Value = "Tropicales y subtropicales"
$ws.Range("I17").Value = 100108003
$ws.Range("J17").Value = "Maracuyá"
$ws.Range("K17").Value = "Sin especificar"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 30
$ws.Range("N17").Value = 32000
$ws.Range("O17").Value = 32000
$ws.Range("P17").Value = 32000
$ws.Range("Q17").Value = "$/caja 18 kilos"
$ws.Range("R17").Value = "Perú"
$ws.Range("S17").Value = 1778
$ws.Range("T17").Value = 18

# Date cells use the workbook's date number format (style index 2, applied to column D)
$dateFormat = $ws.Range("D18").NumberFormat
$ws.Range("D16").NumberFormat = $dateFormat
$ws.Range("D17").NumberFormat = $dateFormat
